$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 57.14035266666667
$ws.Range("H2").Value = 171.421058
$ws.Range("I2").Value = 0.7274038390747541
$ws.Range("J2").Value = 0.7274038390747541
$ws.Range("M2").Value = 15.75563966666667
$ws.Range("N2").Value = 47.266919
$ws.Range("O2").Value = 0.3220556913988901
$ws.Range("P2").Value = 0.32205569139889
$ws.Range("Q2").Value = 900.2828070422559
$ws.Range("R2").Value = 8102.545263380303
$ws.Range("S2").Value = 0.2342645463194269
$ws.Range("T2").Value = 0.2342645463194269
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 57.14035266666667
$ws.Range("H3").Value = 171.421058
$ws.Range("I3").Value = 0.7274038390747541
$ws.Range("J3").Value = 0.7274038390747541
$ws.Range("O3").Value = 0.5509544596378365
$ws.Range("P3").Value = 0.5509544596378364
$ws.Range("Q3").Value = 1540.15234235016
$ws.Range("R3").Value = 13861.37108115144
$ws.Range("S3").Value = 0.4007663890959189
$ws.Range("T3").Value = 0.4007663890959189
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 57.14035266666667
$ws.Range("H4").Value = 171.421058
$ws.Range("I4").Value = 0.7274038390747541
$ws.Range("J4").Value = 0.7274038390747541
$ws.Range("O4").Value = 0.1269898489632735
$ws.Range("P4").Value = 0.1269898489632735
$ws.Range("Q4").Value = 354.9907073336758
$ws.Range("R4").Value = 3194.916366003082
$ws.Range("S4").Value = 0.09237290365940835
$ws.Range("T4").Value = 0.09237290365940834
$ws.Range("I5").Value = 0.08622113322131104
$ws.Range("J5").Value = 0.08622113322131104
$ws.Range("M5").Value = 15.75563966666667
$ws.Range("N5").Value = 47.266919
$ws.Range("O5").Value = 0.3220556913988901
$ws.Range("P5").Value = 0.32205569139889
$ws.Range("Q5").Value = 106.7129422104534
$ws.Range("R5").Value = 960.4164798940809
$ws.Range("S5").Value = 0.02776800667278514
$ws.Range("T5").Value = 0.02776800667278513
$ws.Range("I6").Value = 0.08622113322131104
$ws.Range("J6").Value = 0.08622113322131104
$ws.Range("O6").Value = 0.5509544596378365
$ws.Range("P6").Value = 0.5509544596378364
$ws.Range("S6").Value = 0.04750391786330934
$ws.Range("T6").Value = 0.04750391786330933
$ws.Range("I7").Value = 0.08622113322131104
$ws.Range("J7").Value = 0.08622113322131104
$ws.Range("O7").Value = 0.1269898489632735
$ws.Range("P7").Value = 0.1269898489632735
$ws.Range("S7").Value = 0.01094920868521657
$ws.Range("T7").Value = 0.01094920868521657
$ws.Range("I8").Value = 0.1863750277039348
$ws.Range("J8").Value = 0.1863750277039348
$ws.Range("M8").Value = 15.75563966666667
$ws.Range("N8").Value = 47.266919
$ws.Range("O8").Value = 0.3220556913988901
$ws.Range("P8").Value = 0.32205569139889
$ws.Range("Q8").Value = 230.669985626283
$ws.Range("R8").Value = 2076.029870636547
$ws.Range("S8").Value = 0.060023138406678
$ws.Range("T8").Value = 0.060023138406678
$ws.Range("I9").Value = 0.1863750277039348
$ws.Range("J9").Value = 0.1863750277039348
$ws.Range("O9").Value = 0.5509544596378365
$ws.Range("P9").Value = 0.5509544596378364
$ws.Range("S9").Value = 0.1026841526786082
$ws.Range("T9").Value = 0.1026841526786082
$ws.Range("I10").Value = 0.1863750277039348
$ws.Range("J10").Value = 0.1863750277039348
$ws.Range("O10").Value = 0.1269898489632735
$ws.Range("P10").Value = 0.1269898489632735
$ws.Range("S10").Value = 0.0236677366186486
$ws.Range("T10").Value = 0.0236677366186486
